$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Move to location (7, 5) and remove the screws.'
$ws.Range("B2").Value = 26.214176
$ws.Range("C2").Value = 3773
$ws.Range("D2").Value = "'0.00738"
$ws.Range("E2").Value = '304fa2a7-8596-4062-8bc3-e4730f4005c8'

$ws.Range("A3").Value = 'Move Robot2 to location (11, 8) and remove the toolkit.'
$ws.Range("B3").Value = 24.942843
$ws.Range("C3").Value = 3836
$ws.Range("D3").Value = "'0.0081"
$ws.Range("E3").Value = '21fb56e1-03e9-4472-9547-8bb46a360ad8'

$ws.Range("A4").Value = 'Move Robot26 to location (4, 4) and remove the liquid spill.'
$ws.Range("B4").Value = 22.071861
$ws.Range("C4").Value = 3931
$ws.Range("D4").Value = "'0.0084"
$ws.Range("E4").Value = '5cf21d66-8c2b-4e4f-9299-47008f91dc54'

$ws.Range("A5").Value = 'Move Robot42 to location (9, 1) and remove the large debris.'
$ws.Range("B5").Value = 27.840383
$ws.Range("C5").Value = 4000
$ws.Range("D5").Value = "'0.00867"
$ws.Range("E5").Value = '7eb01efa-8d8c-48c1-9eb9-bf80c7f64ab1'

$ws.Range("A6").Value = 'Move Robot50 to location (7, 11) and remove the dust.'
$ws.Range("B6").Value = 30.961912
$ws.Range("C6").Value = 4584
$ws.Range("D6").Value = "'0.00963"
$ws.Range("E6").Value = '5e73f4cc-a413-40ab-86b9-9f8c9cb39d66'

$ws.Range("A7").Value = 'Move Robot41 to location (6, 12) and remove the grass.'
$ws.Range("B7").Value = 22.119066
$ws.Range("C7").Value = 3801
$ws.Range("D7").Value = "'0.00756"
$ws.Range("E7").Value = '05f759ff-bac4-47bd-a868-5aad31d009e9'

$ws.Range("A8").Value = 'Move Robot50 to location (3, 1) and remove the small debris.'
$ws.Range("B8").Value = 32.655067
$ws.Range("C8").Value = 4648
$ws.Range("D8").Value = "'0.01002"
$ws.Range("E8").Value = '345ccc0b-addb-4eed-b8d8-35a7d8cd3944'

$ws.Range("A9").Value = 'Move Robot13 to location (1, 4) and remove the vehicle.'
$ws.Range("B9").Value = 24.222663
$ws.Range("C9").Value = 3819
$ws.Range("D9").Value = "'0.0081"
$ws.Range("E9").Value = 'f9cf4a02-7c30-4326-95f4-3ae0b0494580'

$ws.Range("A10").Value = 'Move Robot13 to location (11, 1) and remove the construction materials.'
$ws.Range("B10").Value = 28.322895
$ws.Range("C10").Value = 3921
$ws.Range("D10").Value = "'0.00774"
$ws.Range("E10").Value = '7420de6d-4ab9-4dad-ae4b-33290d2527dc'

$ws.Range("A11").Value = 'Move Robot14 to location (2, 10) and remove the tree branches.'
$ws.Range("B11").Value = 27.654019
$ws.Range("C11").Value = 3848
$ws.Range("D11").Value = "'0.00792"
$ws.Range("E11").Value = '6426e022-a3ec-41f9-bc7c-52324826cc9f'

